$d = $word.ActiveDocument

# --- Edit 1: split the "Researching...studied for QC" run into two runs ---
# (Angel's standup note becomes two separate runs even though the visible text
# and formatting stay the same -- mirrors a real edit where the author
# clicked mid-sentence and retyped.)
$text = $d.Content.Text
$needle = "Researching for image feature to chat and studied for QC"
$idx = $text.IndexOf($needle)
if ($idx -ge 0) {
    $splitAt = $idx + ("Researching for image feature to chat and s").Length
    $endAt = $idx + $needle.Length
    $tail = $d.Range($splitAt, $endAt)
    # Toggling a character property and reverting it forces Word to
    # materialize a run boundary at $splitAt without changing the
    # effective formatting of the text.
    $tail.Font.Bold = 1
    $tail.Font.Bold = 0
}

# --- Edit 2: remove the "Kyla/Adam/Alejandro working on..." + "Blockers-"
# block that followed the "Login" heading (sprint_2_python merge revert). ---
$idx = 1
$startPara = $null
$endPara = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -eq "Login`r") {
        $startPara = $idx + 1
    }
    if ($t -eq "Alejandro: clarify dao testing documentation`r") {
        $endPara = $idx
    }
    $idx = $idx + 1
}

if ($startPara -ne $null -and $endPara -ne $null -and $startPara -le $endPara) {
    $rStart = $d.Paragraphs($startPara).Range.Start
    $rEnd = $d.Paragraphs($endPara).Range.End
    $d.Range($rStart, $rEnd).Delete()
}
